# Refresh cryptos list (GitHub Actions scheduled data update).
# Updates Price (D) and Volume(1h) (E) columns for the latest quotes, and
# reflects Aave overtaking TheSandbox in the rank-40/41 slots.
#
# NumberFormat = "@" is applied before writing any Price value that would
# otherwise parse as a plain number (e.g. "1.000", "0.7758") so it is
# stored as literal text, matching the original data which always stores
# Price as text (values like "1.000" must keep their trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.880.14"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.892.34"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7758"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.60"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3145"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07371"
$ws.Range("E9").Value = "  +3.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.37"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08139"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7665"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.483"
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("D14").Value = "1.877.95"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.48"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.227"
$ws.Range("E16").Value = "  +4.96%  "
$ws.Range("D17").Value = "29.869.38"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.98"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.88"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007860"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.155"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "2.129.04"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1579"
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.448"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.64"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.82"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.046"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.448"
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.548"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.508"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05604"
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.100"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.253"
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7641"
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.000"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.646"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01936"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.791"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "1.159.25"
$ws.Range("E41").Value = "  +12.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.977"
$ws.Range("E44").Value = "  +1.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8547"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.908"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.157"
$ws.Range("E48").Value = "  +6.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.00"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.549"
$ws.Range("E51").Value = "  +0.65%  "

# Rows 42 and 43: Aave and TheSandbox swap positions with updated data
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.48"
$ws.Range("E42").Value = "  +2.96%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4470"
$ws.Range("E43").Value = "  +0.21%  "
